# 1. Merge "Most people won't appear o" + bookmark + "n any list due to the
#    100-page max." into a single clean run. Matching (and replacing with)
#    the full combined sentence makes Find span across the old mid-sentence
#    bookmark, which collapses the two runs into one and drops the old
#    bookmark (it gets re-created later, at the end of the new paragraph).
$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Most people won" + [char]0x2019 + "t appear on any list due to the 100-page max.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Most people won" + [char]0x2019 + "t appear on any list due to the 100-page max.",
    2) | Out-Null

# 2. Drop the stray trailing space at the end of the "We can have some
#    gauge..." paragraph.
$d.Content.Find.Execute(
    "quite strong and reliable. ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "quite strong and reliable.",
    2) | Out-Null

# 3. Insert the brand-new trailing bullet paragraph right after the
#    "We can have some gauge..." paragraph, inheriting its ListParagraph /
#    ilvl=2 / numId=5 formatting.
$hostPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*quite strong and reliable.*") {
        $hostPara = $cand
    }
}
$hostRange = $hostPara.Range
$hostRange.InsertParagraphAfter()

$newParaIndex = $hostPara.Index + 1
$newPara = $d.Paragraphs.Item($newParaIndex)
$newRange = $newPara.Range
# Trailing sentinel character "Z" so we can anchor the "_GoBack" bookmark at
# the true end of the sentence (collapsed ranges sitting exactly at a
# paragraph's final text position don't host new bookmarks reliably in this
# engine), then trim the sentinel back off afterwards.
$newRange.Text = "Github does have a so-called abuse prevention mechanism in place which occasionally prevented my scraper from obtaining results during a scraping run. It" + [char]0x2019 + "s not obvious this matters, but worth mentioning.Z"

$f = $d.Content
$found3 = $f.Find.Execute("but worth mentioning.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $f.End

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$sentinelRange = $d.Range($endPos, $endPos + 1)
$sentinelRange.Text = ""

